# Update the "as_of_utc" refresh timestamp (column AA) for every data row
# on the "Главные" and "Линейные" sheets, reflecting the newer publish run.
$wb = $excel.ActiveWorkbook

$oldTimestamp = "2025-12-13 03:02:44"
$newTimestamp = "2025-12-13 07:02:53"

$sheetNames = @("Главные", "Линейные")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    if ($lastRow -lt 2) { $lastRow = 2 }

    for ($row = 2; $row -le $lastRow; $row++) {
        $cell = $ws.Cells.Item($row, 27)   # column AA
        if ($cell.Value2 -eq $oldTimestamp) {
            $cell.Value = $newTimestamp
        }
    }
}
